# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet before "ODI Batting" with the
#    player's ID/NAME/BATTING_HAND/BOWL_STYLE.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, and replace the full
#    howstat.com scorecard URL values with the bare numeric match code.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 1. New "Player Info" sheet, placed before "ODI Batting" ---------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "4438"
$playerInfo.Range("B2").Value = "Gavin T Main"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# Re-resolve the batting/bowling sheets by name: inserting a new sheet shifts
# worksheet indices, and stale references can end up pointing at the sheet
# that now occupies the old index.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D) ------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        $cell.Value = $code
    }
}

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B) ------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '^.*MatchCode=', ''
        $cell.Value = $code
    }
}
